$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 3539.4546
$ws.Range("I31").Value = 2793.5
$ws.Range("K31").Value = 8380.5
$ws.Range("M31").Value = -8150.5

$ws.Range("H33").Value = 175.5
$ws.Range("I33").Value = 183.14285
$ws.Range("K33").Value = 183.14285
$ws.Range("M33").Value = 45.85714999999999

$ws.Range("H55").Value = 795.3333
$ws.Range("J55").Value = 795.3333
$ws.Range("L55").Value = 795.3333
$ws.Range("N55").Value = -1223.3333

$ws.Range("H96").Value = 997.5714
$ws.Range("I96").Value = 972.5
$ws.Range("K96").Value = 2917.5
$ws.Range("M96").Value = -1544.5

$ws.Range("H98").Value = 989.7895
$ws.Range("I98").Value = 947.9231
$ws.Range("K98").Value = 947.9231
$ws.Range("M98").Value = 550.0769

$ws.Range("H111").Value = 2228.6667
$ws.Range("I111").Value = 1666
$ws.Range("J111").Value = 2791.3333
$ws.Range("K111").Value = 4998
$ws.Range("L111").Value = 8373.999899999999
$ws.Range("M111").Value = -1931
$ws.Range("N111").Value = -14507.9999

$ws.Range("H122").Value = 989.7895
$ws.Range("I122").Value = 947.9231
$ws.Range("K122").Value = 2843.7693
$ws.Range("M122").Value = -393.7692999999999

$ws.Range("H132").Value = 4024.0334
$ws.Range("I132").Value = 3738.7036
$ws.Range("K132").Value = 11216.1108
$ws.Range("M132").Value = -8686.110799999999

$ws.Range("H138").Value = 4469.6284
$ws.Range("J138").Value = 4541.278
$ws.Range("L138").Value = 13623.834
$ws.Range("N138").Value = -23903.834

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1273.9054
$ws.Range("I32").Value = 1085.8767
$ws.Range("K32").Value = 1085.8767
$ws.Range("M32").Value = -798.8767

$ws.Range("H61").Value = 5672.212
$ws.Range("I61").Value = 2546.96
$ws.Range("K61").Value = 2546.96
$ws.Range("M61").Value = -2334.96

$ws.Range("H63").Value = 3218.25
$ws.Range("I63").Value = 2492.5
$ws.Range("K63").Value = 2492.5
$ws.Range("M63").Value = -1806.5

$ws.Range("H66").Value = 3218.25
$ws.Range("I66").Value = 2492.5
$ws.Range("K66").Value = 12462.5
$ws.Range("M66").Value = -9030.5

$ws.Range("H74").Value = 1947.9592
$ws.Range("I74").Value = 1590.2106
$ws.Range("J74").Value = 2174.5334
$ws.Range("K74").Value = 1590.2106
$ws.Range("L74").Value = 2174.5334
$ws.Range("M74").Value = -716.2106000000001
$ws.Range("N74").Value = -3922.5334

$ws.Range("H77").Value = 1947.9592
$ws.Range("I77").Value = 1590.2106
$ws.Range("J77").Value = 2174.5334
$ws.Range("K77").Value = 7951.053000000001
$ws.Range("L77").Value = 10872.667
$ws.Range("M77").Value = -3583.053000000001
$ws.Range("N77").Value = -19608.667

$ws.Range("H136").Value = 5672.212
$ws.Range("I136").Value = 2546.96
$ws.Range("K136").Value = 7640.88
$ws.Range("M136").Value = -5090.88

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 17784.705
$ws.Range("I94").Value = 8973
$ws.Range("K94").Value = 8973
$ws.Range("M94").Value = -8522

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2372.625
$ws.Range("J31").Value = 2999.6667
$ws.Range("L31").Value = 2999.6667
$ws.Range("N31").Value = -3589.6667

$ws.Range("H34").Value = 2372.625
$ws.Range("J34").Value = 2999.6667
$ws.Range("L34").Value = 2999.6667
$ws.Range("N34").Value = -3403.6667

$ws.Range("H58").Value = 3335.5
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3335.5
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3335.5
$ws.Range("N58").Value = -3741.5
$ws.Range("M58").ClearContents()

$ws.Range("H99").Value = 1999.3334
$ws.Range("I99").Value = 1899.5
$ws.Range("K99").Value = 1899.5
$ws.Range("M99").Value = -401.5

$ws.Range("H105").Value = 1628.5
$ws.Range("I105").Value = 1628.5
$ws.Range("K105").Value = 1628.5
$ws.Range("M105").Value = 118.5

$ws.Range("H122").Value = 6161.3125
$ws.Range("I122").Value = 5053
$ws.Range("K122").Value = 15159
$ws.Range("M122").Value = -12709

$ws.Range("H126").Value = 1999.3334
$ws.Range("I126").Value = 1899.5
$ws.Range("K126").Value = 5698.5
$ws.Range("M126").Value = -3228.5

$ws.Range("H132").Value = 2137.6667
$ws.Range("I132").Value = 2137.6667
$ws.Range("K132").Value = 6413.000100000001
$ws.Range("M132").Value = -3883.000100000001

$ws.Range("H136").Value = 3335.5
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3335.5
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 10006.5
$ws.Range("N136").Value = -15106.5
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 537773.75
$ws.Range("I4").Value = 720797
$ws.Range("J4").Value = 145581
$ws.Range("K4").Value = 2162391
$ws.Range("L4").Value = 436743
$ws.Range("M4").Value = -2162279
$ws.Range("N4").Value = -436967

$ws.Range("H107").Value = 819.7
$ws.Range("I107").Value = 248.6
$ws.Range("J107").Value = 1390.8
$ws.Range("K107").Value = 745.8
$ws.Range("L107").Value = 4172.4
$ws.Range("M107").Value = 1174.2
$ws.Range("N107").Value = -8012.4

$ws.Range("H109").Value = 3369.1538
$ws.Range("I109").Value = 2918.0908
$ws.Range("K109").Value = 8754.2724
$ws.Range("M109").Value = -7714.2724

$ws.Range("H121").Value = 95801.21000000001
$ws.Range("J121").Value = 95478.234
$ws.Range("L121").Value = 286434.702
$ws.Range("N121").Value = -289054.702

$ws.Range("H126").Value = 9006.25
$ws.Range("I126").Value = 9006.25
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 27018.75
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -22078.75
$ws.Range("N126").ClearContents()

$ws.Range("H129").Value = 2199.25
$ws.Range("I129").Value = 882.6
$ws.Range("K129").Value = 2647.8
$ws.Range("M129").Value = 2352.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1144.5
$ws.Range("I107").Value = 144
$ws.Range("J107").Value = 1478
$ws.Range("K107").Value = 144
$ws.Range("L107").Value = 1478
$ws.Range("M107").Value = 1776
$ws.Range("N107").Value = -5318

$ws.Range("H113").Value = 4583.8887
$ws.Range("I113").Value = 1353.5
$ws.Range("J113").Value = 5506.857
$ws.Range("K113").Value = 1353.5
$ws.Range("L113").Value = 5506.857
$ws.Range("M113").Value = 816.5
$ws.Range("N113").Value = -9846.857

$ws.Range("H122").Value = 2870.8462
$ws.Range("I122").Value = 2557.9412
$ws.Range("J122").Value = 3461.889
$ws.Range("K122").Value = 7673.823600000001
$ws.Range("L122").Value = 10385.667
$ws.Range("M122").Value = -5223.823600000001
$ws.Range("N122").Value = -15285.667

$ws.Range("H138").Value = 74999.89999999999
$ws.Range("J138").Value = 74999.89999999999
$ws.Range("L138").Value = 74999.89999999999
$ws.Range("N138").Value = -85279.89999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2093.0833
$ws.Range("I22").Value = 1131
$ws.Range("J22").Value = 3440
$ws.Range("K22").Value = 1131
$ws.Range("L22").Value = 3440
$ws.Range("M22").Value = -836
$ws.Range("N22").Value = -4030

$ws.Range("H27").Value = 2093.0833
$ws.Range("I27").Value = 1131
$ws.Range("J27").Value = 3440
$ws.Range("K27").Value = 1131
$ws.Range("L27").Value = 3440
$ws.Range("M27").Value = -1024
$ws.Range("N27").Value = -3654

$ws.Range("H55").Value = 1698.625
$ws.Range("J55").Value = 1798.4286
$ws.Range("L55").Value = 1798.4286
$ws.Range("N55").Value = -2144.4286

$ws.Range("H61").Value = 2679.4
$ws.Range("I61").Value = 2649.25
$ws.Range("J61").Value = 2800
$ws.Range("K61").Value = 2649.25
$ws.Range("L61").Value = 2800
$ws.Range("M61").Value = -2447.25
$ws.Range("N61").Value = -3204

$ws.Range("H113").Value = 2679.4
$ws.Range("I113").Value = 2649.25
$ws.Range("J113").Value = 2800
$ws.Range("K113").Value = 2649.25
$ws.Range("L113").Value = 2800
$ws.Range("M113").Value = -479.25
$ws.Range("N113").Value = -7140

$ws.Range("H122").Value = 4225.7803
$ws.Range("I122").Value = 2941.5557
$ws.Range("K122").Value = 8824.667099999999
$ws.Range("M122").Value = -6374.667099999999

$ws.Range("H132").Value = 1954.6552
$ws.Range("I132").Value = 1830.6857
$ws.Range("K132").Value = 5492.0571
$ws.Range("M132").Value = -2962.0571

$ws.Range("H134").Value = 82111.11
$ws.Range("J134").Value = 82111.11
$ws.Range("L134").Value = 82111.11
$ws.Range("N134").Value = -92251.11

$ws.Range("H138").Value = 75000
$ws.Range("J138").Value = 75000
$ws.Range("L138").Value = 75000
$ws.Range("N138").Value = -85280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1511.7894
$ws.Range("I113").Value = 1444.4546
$ws.Range("J113").Value = 1604.375
$ws.Range("K113").Value = 4333.3638
$ws.Range("L113").Value = 4813.125
$ws.Range("M113").Value = -2163.3638
$ws.Range("N113").Value = -9153.125

$ws.Range("H132").Value = 4438.921
$ws.Range("I132").Value = 4421.2188
$ws.Range("K132").Value = 13263.6564
$ws.Range("M132").Value = -10733.6564
